{"js": "// Append the \"Business Logic constraints:\" section to the end of the\n// document body, after the existing final (empty) paragraph.\n//\n// Target shape (from the OOXML diff):\n//   - 5 empty paragraphs carrying the same bold/blue-gray/36pt paragraph\n//     mark formatting as the document's existing trailing paragraph.\n//   - A \"Business Logic constraints:\" heading (bold, 40 half-pt-ish /\n//     20pt, i.e. w:sz 40) made of 3 runs, the first one carrying a\n//     <w:lastRenderedPageBreak/> marker.\n//   - A \"Ronald Around (Owner)\" sub-heading (bold, 16pt / w:sz 32).\n//   - Two 12pt (w:sz 24) body paragraphs describing the constraints; the\n//     last one's paragraph mark (w:pPr/w:rPr) is bold + colored even\n//     though its own runs are not (matches the source document exactly).\n//\n// We build this as a minimal \"flat OPC\" WordprocessingML package and feed\n// it to Body.insertOoxml(), which lets us control every run/paragraph\n// property exactly instead of relying on inherited paragraph-mark\n// formatting from higher level APIs (Office.js font.* setters would\n// otherwise stamp identical formatting on both the paragraph mark and the\n// runs, which is not what the source document has for the last\n// paragraph).\n\nconst emptySpacerParagraphXml =\n  '<w:p><w:pPr><w:rPr><w:b/><w:bCs/>' +\n  '<w:color w:val=\"262626\" w:themeColor=\"text1\" w:themeTint=\"D9\"/>' +\n  '<w:sz w:val=\"36\"/><w:szCs w:val=\"36\"/></w:rPr></w:pPr></w:p>';\n\nconst headingRunXml = (text, extra) =>\n  '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"40\"/><w:szCs w:val=\"40\"/></w:rPr>' +\n  (extra || '') +\n  '<w:t xml:space=\"preserve\">' + text + '</w:t></w:r>';\n\nconst headingParagraphXml =\n  '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val=\"40\"/><w:szCs w:val=\"40\"/></w:rPr></w:pPr>' +\n  headingRunXml('Business ', '<w:lastRenderedPageBreak/>') +\n  headingRunXml('Logic ') +\n  '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"40\"/><w:szCs w:val=\"40\"/></w:rPr><w:t>constraints:</w:t></w:r>' +\n  '</w:p>';\n\nconst ownerParagraphXml =\n  '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/></w:rPr>' +\n  '<w:t>Ronald Around (Owner)</w:t></w:r></w:p>';\n\nconst soldPriceParagraphXml =\n  '<w:p><w:pPr><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr>' +\n  '<w:t>Only Ronald Around (Owner) can enter sold prices that are less than or equal to 95% of the invoice price.</w:t></w:r></w:p>';\n\nconst laborChargeParagraphXml =\n  '<w:p><w:pPr><w:rPr><w:b/><w:bCs/>' +\n  '<w:color w:val=\"262626\" w:themeColor=\"text1\" w:themeTint=\"D9\"/>' +\n  '<w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr>' +\n  '<w:t>Only Ronald Around (Owner) can update the labor charges on a repair to a value less than their previous value</w:t></w:r>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r>' +\n  '</w:p>';\n\nconst newParagraphsXml =\n  emptySpacerParagraphXml.repeat(5) +\n  headingParagraphXml +\n  ownerParagraphXml +\n  soldPriceParagraphXml +\n  laborChargeParagraphXml;\n\nconst flatOpcXml =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<?mso-application progid=\"Word.Document\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  newParagraphsXml +\n  '<w:sectPr><w:pgSz w:w=\"12240\" w:h=\"15840\"/></w:sectPr>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\n// Insert at the very end of the body, after the existing (untouched)\n// trailing paragraph.\ncontext.document.body.insertOoxml(flatOpcXml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Append the \"Business Logic constraints:\" section to the end of the\n# document body, after the existing final (empty) paragraph.\n#\n# Target shape (from the OOXML diff):\n#   - 5 empty paragraphs carrying the same bold/blue-gray/36pt paragraph\n#     mark formatting as the document's existing trailing paragraph.\n#   - A \"Business Logic constraints:\" heading (bold, w:sz 40) made of 3\n#     runs, the first one carrying a <w:lastRenderedPageBreak/> marker.\n#   - A \"Ronald Around (Owner)\" sub-heading (bold, w:sz 32).\n#   - Two w:sz 24 body paragraphs describing the constraints; the last\n#     one's paragraph mark (w:pPr/w:rPr) is bold + colored even though its\n#     own runs are not (matches the source document exactly).\n#\n# We build this as a minimal \"flat OPC\" WordprocessingML package and feed\n# it to Range.InsertXML(), which gives exact control over every\n# run/paragraph property instead of relying on Word inheriting the\n# adjacent paragraph mark formatting onto new runs.\n\n$d = $word.ActiveDocument\n\n$emptySpacerParagraphXml = '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:color w:val=\"262626\" w:themeColor=\"text1\" w:themeTint=\"D9\"/><w:sz w:val=\"36\"/><w:szCs w:val=\"36\"/></w:rPr></w:pPr></w:p>'\n\n$headingParagraphXml = '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val=\"40\"/><w:szCs w:val=\"40\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"40\"/><w:szCs w:val=\"40\"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">Business </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"40\"/><w:szCs w:val=\"40\"/></w:rPr><w:t xml:space=\"preserve\">Logic </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"40\"/><w:szCs w:val=\"40\"/></w:rPr><w:t>constraints:</w:t></w:r></w:p>'\n\n$ownerParagraphXml = '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/></w:rPr><w:t>Ronald Around (Owner)</w:t></w:r></w:p>'\n\n$soldPriceParagraphXml = '<w:p><w:pPr><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>Only Ronald Around (Owner) can enter sold prices that are less than or equal to 95% of the invoice price.</w:t></w:r></w:p>'\n\n$laborChargeParagraphXml = '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:color w:val=\"262626\" w:themeColor=\"text1\" w:themeTint=\"D9\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>Only Ronald Around (Owner) can update the labor charges on a repair to a value less than their previous value</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r></w:p>'\n\n$newParagraphsXml = ($emptySpacerParagraphXml * 5) + $headingParagraphXml + $ownerParagraphXml + $soldPriceParagraphXml + $laborChargeParagraphXml\n\n$flatOpcXml = @'\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<?mso-application progid=\"Word.Document\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n__NEW_PARAGRAPHS__\n<w:sectPr><w:pgSz w:w=\"12240\" w:h=\"15840\"/></w:sectPr>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n'@\n\n$flatOpcXml = $flatOpcXml.Replace('__NEW_PARAGRAPHS__', $newParagraphsXml)\n\n# Create a fresh landing paragraph right after the existing trailing\n# paragraph, leaving that paragraph completely untouched, then replace the\n# (empty) landing paragraph's content with our XML block. Word always\n# merges the *last* paragraph mark of inserted XML into whatever paragraph\n# mark was at the insertion point, so inserting at the very start of this\n# disposable empty paragraph makes our last new paragraph take over its\n# slot cleanly, while every paragraph before it becomes fully new.\n$lastPara = $d.Paragraphs.Last\n$tailRange = $lastPara.Range\n$tailRange.Collapse(0)\n$tailRange.InsertParagraphAfter()\n\n$landingPara = $d.Paragraphs.Last\n$landingRange = $landingPara.Range\n$landingRange.Collapse(1)\n$landingRange.InsertXML($flatOpcXml)\n"}
